$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("D2:E51")
$rng.NumberFormat = "@"

$ws.Range("D2").Value = "29.197.42"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").Value = "1.834.90"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("D4").Value = "0.9989"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "241.95"
$ws.Range("E5").Value = "  +0.78%  "
$ws.Range("D6").Value = "0.6653"
$ws.Range("E6").Value = "  -2.31%  "
$ws.Range("D7").Value = "0.9998"
$ws.Range("D8").Value = "0.07415"
$ws.Range("E8").Value = "  -0.63%  "
$ws.Range("D9").Value = "0.2933"
$ws.Range("E9").Value = "  -1.74%  "
$ws.Range("E10").Value = "  -0.87%  "
$ws.Range("D12").Value = "1.836.03"
$ws.Range("E12").Value = "  +1.45%  "
$ws.Range("D13").Value = "5.000"
$ws.Range("E13").Value = "  -0.30%  "
$ws.Range("D14").Value = "0.6680"
$ws.Range("E14").Value = "  -1.31%  "
$ws.Range("D15").Value = "83.03"
$ws.Range("E15").Value = "  -4.09%  "
$ws.Range("E16").Value = "  -0.63%  "
$ws.Range("D17").Value = "0.000008361"
$ws.Range("E17").Value = "  +1.10%  "
$ws.Range("D18").Value = "29.194.03"
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("D19").Value = "2.086.23"
$ws.Range("E19").Value = "  +2.14%  "
$ws.Range("D20").Value = "228.16"
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("D21").Value = "12.46"
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("E23").Value = "  -2.84%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").Value = "159.10"
$ws.Range("E25").Value = "  -1.12%  "
$ws.Range("E26").Value = "  -2.26%  "
$ws.Range("D27").Value = "8.616"
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("D29").Value = "1.512"
$ws.Range("E29").Value = "  +0.70%  "
$ws.Range("D30").Value = "4.113"
$ws.Range("E30").Value = "  -3.21%  "
$ws.Range("D31").Value = "4.042"
$ws.Range("E32").Value = "  -0.36%  "
$ws.Range("D33").Value = "0.05284"
$ws.Range("E33").Value = "  -2.32%  "
$ws.Range("D34").Value = "1.865"
$ws.Range("E34").Value = "  +0.53%  "
$ws.Range("D35").Value = "0.7461"
$ws.Range("E35").Value = "  -0.85%  "
$ws.Range("E36").Value = "  +0.80%  "
$ws.Range("D37").Value = "2.649"
$ws.Range("E37").Value = "  -1.19%  "
$ws.Range("D38").Value = "1.292.31"
$ws.Range("E38").Value = "  -0.94%  "
$ws.Range("E39").Value = "  -0.79%  "
$ws.Range("D40").Value = "2.737"
$ws.Range("E40").Value = "  +0.76%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").Value = "5.876"
$ws.Range("E42").Value = "  -2.90%  "
$ws.Range("D43").Value = "0.08364"
$ws.Range("E43").Value = "  -1.42%  "
$ws.Range("D44").Value = "0.9991"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").Value = "102.24"
$ws.Range("E45").Value = "  -2.53%  "
$ws.Range("D46").Value = "1.985.60"
$ws.Range("E46").Value = "  +1.27%  "
$ws.Range("D48").Value = "1.758"
$ws.Range("E49").Value = "  -1.09%  "
$ws.Range("D50").Value = "62.98"
$ws.Range("E50").Value = "  -1.02%  "
$ws.Range("D51").Value = "0.05873"
$ws.Range("E51").Value = "  -0.88%  "

$rng.ClearFormats()
